$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: every value in this sheet (including numbers) is stored as text
# (originally t="inlineStr" cells). A plain .Value assignment of a
# numeric-looking string gets auto-coerced into a real number by Excel,
# which would e.g. strip the leading zeros from SKU codes. So for every
# cell we write, the NumberFormat is set to Text ("@") first to keep the
# written value as text, matching the rest of the sheet.

# --- Update existing rows to reflect reduced quantities/totals ---

# Row 2: Natalie's - Orange Juice : Quantity 4 -> 2, Total 95.00 -> 47.50
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "47.50"

# Row 4: Natalie's - Orange Pineapple : Quantity 3 -> 2, Total 39.00 -> 26.00
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "26.00"

# Row 5: Natalie's - Orange Mango : Quantity 3 -> 2, Total 39.00 -> 26.00
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "2"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "26.00"

# Row 7: Natalie's - Honey Tangerine : Quantity 2 -> 1, Total 28.00 -> 14.00
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "1"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "14.00"

# --- Append two new order line items (rows 9 & 10) ---

# Row 9: Natalie's - Lemonade
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "004011"
$ws.Range("B9").Value = "Natalie's - Lemonade"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "1"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "9.25"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "9.25"

# Row 10: Natalie's - Strawberry Lemonade
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "004014"
$ws.Range("B10").Value = "Natalie's - Strawberry Lemonade"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "1"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "13.90"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "13.90"
